# FHIR-34724: update to lab and procedures
# - Patient "deceased" search param renamed to "death-date" (sps + sp_combos sheets)
# - New RelatedPerson.name search param row added to the "sps" sheet
# - RelatedPerson.patient example GET updated to use the actual search param name
# - defined name / sheet view bookkeeping refreshed to match

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "sps" worksheet
# ------------------------------------------------------------------
$sps = $wb.Worksheets.Item("sps")

# 1a. Patient.deceased -> Patient.death-date (row 25). Dependent formulas
#     (L25, AB25) recompute automatically since they reference C25.
$sps.Range("C25").Value2 = "death-date"

# 1b. RelatedPerson.patient (row 113) sample GET now parameterised off C113
#     instead of being hard-coded to "patient".
$sps.Range("Z113").Formula = '="GET [base]/"&B113&"?"&C113&"=1032702"'

# 1c. Insert a brand-new row for RelatedPerson.name above the old row 114
#     (QuestionnaireResponse._id, ... shift down to 115.. automatically,
#     carrying their formulas/values with them unchanged).
$sps.Rows.Item(114).Insert()

$sps.Range("A114").Value2 = 26
$sps.Range("B114").Value2 = "RelatedPerson"
$sps.Range("C114").Value2 = "name"
$sps.Range("D114").Value2 = "SHOULD"
$sps.Range("E114").Value2 = $true
$sps.Range("G114").Formula = '="http://hl7.org/fhir/us/core/StructureDefinition/us-core-"&LOWER(B114)'
$sps.Range("H114").Value2 = "Y"
$sps.Range("I114").Value2 = "Y"
$sps.Range("J114").Value2 = "Y"
$sps.Range("K114").Value2 = "string"
$sps.Range("L114").Formula = '=B114&"."&C114'
$sps.Range("M114").Value2 = "Y"
$sps.Range("O114").Value2 = "Y"

# match formatting (font style "4"/"10") used by the neighbouring Y/Z/AA cells
$sps.Range("Z113").Copy()
$sps.Range("Y114").PasteSpecial(-4122)
$sps.Range("AA113").Copy()
$sps.Range("AA114").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$sps.Range("Y114").Value2 = "support searching for a patient by a server defined search that matches any of the string fields in the HumanName, including family, give, prefix, suffix, suffix, and/or text"
$sps.Range("Z114").Formula = '="GET [base]/"&B114&"?"&C114&"=Mary Shaw"'
$sps.Range("AA114").Formula = '="Fetches a bundle of all "&B114&" resources matching the name"'
$sps.Range("AB114").Formula = '="SearchParameter-us-core-"&LOWER((B114)&"-"&SUBSTITUTE(C114,"_","")&".html")'

# row insert leaves a stray formatted-but-empty cell in column F; drop the format
$sps.Range("F114").ClearFormats()

# restore the autofilter disturbed by the row insert
$sps.AutoFilterMode = $false
$sps.Range("A1:AB120").AutoFilter()

# ------------------------------------------------------------------
# 2. "sp_combos" worksheet - same deceased -> death-date rename
# ------------------------------------------------------------------
$combos = $wb.Worksheets.Item("sp_combos")
$combos.Range("D32").Value2 = "death-date,family"
$combos.Range("J32").Value2 = "GET [base]/Patient?family=Shaw&death-date=2022-07-22"
$combos.Range("J14").Select()

# ------------------------------------------------------------------
# 3. workbook-level defined name tracking the (now one row taller) sps
#    autofilter range
# ------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "sps!_FilterDatabase") {
        $n.RefersTo = "=sps!`$A`$1:`$AB`$120"
    }
}

# ------------------------------------------------------------------
# 4. leave focus back on "sps" (the originally active/selected tab)
#    with its own selection restored
# ------------------------------------------------------------------
$sps.Activate()
$sps.Range("T84").Select()
